$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Angular 2+, Express, Flask, ASP.NET",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Angular 20+, Express, Flask, ASP.NET",
    2)

$d.Content.Find.Execute(
    "Tested With SOAP UI",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Tested with SOAP UI to find more data on the performance of the APIs",
    2)

$d.Content.Find.Execute(
    "Triaged issues with Splunk",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Triaged issues with Splunk. Monitored and triaged apis to identify errors.",
    2)

$d.Content.Find.Execute(
    "Attended a meeting, communicated with the team on Microsoft Teams.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Communicated regularly through meetings on Teams to review findings",
    2)

$d.Content.Find.Execute(
    "Trained in Spring Boot, Angular, Amazon Web Services",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Trained in Spring Boot, Angular, and Amazon Web Services",
    2)
